# repull data, push all data, mean calculation
# Update the "dSF" column (F) values for each data row on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of row number -> new value for column F ("dSF")
$updates = @{
    2  = -3
    3  = -1
    4  = -4
    5  = 3
    6  = 1
    7  = -2
    8  = 4
    9  = -2
    10 = -2
    11 = 1
    12 = -1
    14 = 0
    15 = 2
    16 = 1
    17 = 2
    18 = -2
    19 = -1
    20 = -6
    21 = -1
    23 = -2
    24 = -9
    25 = -4
    27 = -1
    29 = 1
    30 = 2
    31 = -1
    32 = -2
    33 = -1
    35 = 0
    36 = -5
    37 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
